$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row for "شامبو كلير للرجال 180مل" (row 79) -
# subsequent rows shift up, merged cells adjust accordingly.
$ws.Rows("79:79").Delete()

# Column A holds the sequential item number ("م") which always equals
# (row - 6); the row shift above also advanced these numbers by one, so
# restore the original running sequence (73, 74, 75, 76) for the rows
# that moved up.
$ws.Range("A79").Value = 73
$ws.Range("A80").Value = 74
$ws.Range("A81").Value = 75
$ws.Range("A82").Value = 76

# The totals row (now row 83) needs to reflect the removed item's price
# (80.00) being subtracted from the previous total of 6064.915.
$ws.Range("P83").Value = 5984.915

# Update the printed timestamp (now row 84, column A) to the new save time.
$ws.Range("A84").Value = "Sunday, 28 September, 2025 8:39 PM"

# Row heights are fixed per row position rather than following the
# shifted content; restore the original per-row heights so they keep
# matching the same pattern as before the deletion.
$ws.Rows("79:79").RowHeight = 25.5
$ws.Rows("80:80").RowHeight = 24.75
$ws.Rows("81:81").RowHeight = 25.5
$ws.Rows("82:82").RowHeight = 25.5
$ws.Rows("83:83").RowHeight = 24.75
$ws.Rows("84:84").RowHeight = 16.5
